# Update the date heading.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-02-13 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-02-14 Friday", 2)

# Update the division-problem answers in the single table, cell by cell
# (row, col, newText) so that duplicate source strings are handled correctly.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "74÷5=14, 4"
$t.Cell(1,2).Range.Text = "44÷6=7, 2"
$t.Cell(1,3).Range.Text = "26÷8=3, 2"
$t.Cell(1,4).Range.Text = "99÷6=16, 3"
$t.Cell(1,5).Range.Text = "11÷8=1, 3"

$t.Cell(5,1).Range.Text = "36÷2=18, 0"
$t.Cell(5,2).Range.Text = "11÷7=1, 4"
$t.Cell(5,3).Range.Text = "86÷2=43, 0"
$t.Cell(5,4).Range.Text = "64÷5=12, 4"
$t.Cell(5,5).Range.Text = "95÷5=19, 0"

$t.Cell(9,1).Range.Text = "25÷9=2, 7"
$t.Cell(9,2).Range.Text = "71÷3=23, 2"
$t.Cell(9,3).Range.Text = "44÷3=14, 2"
$t.Cell(9,4).Range.Text = "17÷7=2, 3"
$t.Cell(9,5).Range.Text = "62÷2=31, 0"

$t.Cell(13,1).Range.Text = "76÷6=12, 4"
$t.Cell(13,2).Range.Text = "85÷7=12, 1"
$t.Cell(13,3).Range.Text = "23÷6=3, 5"
$t.Cell(13,4).Range.Text = "41÷4=10, 1"
$t.Cell(13,5).Range.Text = "57÷4=14, 1"

$t.Cell(17,1).Range.Text = "51÷7=7, 2"
$t.Cell(17,2).Range.Text = "97÷9=10, 7"
$t.Cell(17,3).Range.Text = "70÷9=7, 7"
$t.Cell(17,4).Range.Text = "49÷7=7, 0"
$t.Cell(17,5).Range.Text = "19÷4=4, 3"
